# RN.xlsx edit: fill in blank grade cells with 0, add hyperlinks for
# five more student e-mails (with matching "Hipervinculo" cell styling),
# and move the active selection to A18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegistroNotas")

# ---------------------------------------------------------------------
# 1) Blank numeric grade cells (columns C:I) become explicit 0 values.
# ---------------------------------------------------------------------
$zeroCells = @(
    "G2", "H2", "I2",
    "H3",
    "E4", "F4", "G4", "H4", "I4",
    "H7",
    "E8", "F8", "H8",
    "E9", "F9",
    "E15", "F15",
    "E18",
    "F19",
    "F23",
    "E24", "F24", "H24", "I24"
)
foreach ($addr in $zeroCells) {
    $ws.Range($addr).Value = 0
}

# ---------------------------------------------------------------------
# 2) Add hyperlinks (mailto:) for five more students, matching the
#    pattern already used for the other rows, and apply the workbook's
#    existing "Hipervinculo" cell style (preserving each row's banding
#    fill) so the cells look like the other hyperlinked names.
# ---------------------------------------------------------------------
$newLinks = @{
    "A21" = "mailto:dora.solaresalarcon@postgrado.univalle.edu"
    "A14" = "mailto:ramiro.ocanafernandez@postgrado.univalle.edu"
    "A16" = "mailto:marcos.perezhuanca@postgrado.univalle.edu"
    "A23" = "mailto:guisela.vasquezyanez@postgrado.univalle.edu"
    "A18" = "mailto:cesar.rochacruz@postgrado.univalle.edu"
}
$fillBandedColor = 16246759  # BGR for RGB(E7,E7,F7) - the alternating-row band fill

foreach ($addr in @("A21", "A14", "A16", "A23", "A18")) {
    $target = $newLinks[$addr]
    $cell = $ws.Range($addr)

    $ws.Hyperlinks.Add($cell, $target)

    $cell.Style = "Hipervínculo"
    if ($addr -ne "A18") {
        $cell.Interior.Color = $fillBandedColor
    }
}

# ---------------------------------------------------------------------
# 3) Move the active selection to A18 (matches the saved sheet view).
# ---------------------------------------------------------------------
$ws.Range("A18").Select()
